# Fix Subscriptions Report: remove the country.
#
# "Customer Country" is the header in column F of the "Data" sheet's header
# row. Deleting the whole column removes that header (and its shared string)
# and shifts every later column one position to the left, which is exactly
# what the target diff shows (dimension/autoFilter/_FilterDatabase all shrink
# by one column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Delete column F ("Customer Country"); everything to the right shifts left.
$ws.Columns.Item(6).Delete()

# The autofilter range needs to be re-applied so it reflects the new last
# column (was I1:AA1, now H1:Z1 after the column removal).
$ws.AutoFilterMode = $false
$ws.Range("H1:Z1").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the new
# autofilter range as well.
foreach ($n in $ws.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$H`$1:`$Z`$1"
    }
}
